$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "70.721.31"
$ws.Cells.Item(2, 5).Value = "  -1.67%  "

$ws.Cells.Item(3, 4).Value = "2.569.70"
$ws.Cells.Item(3, 5).Value = "  -4.25%  "

$ws.Cells.Item(4, 5).Value = "  +0.03%  "

$ws.Cells.Item(5, 4).Value = "577.68"
$ws.Cells.Item(5, 5).Value = "  -3.27%  "

$ws.Cells.Item(6, 4).Value = "170.80"
$ws.Cells.Item(6, 5).Value = "  -1.79%  "

$ws.Cells.Item(7, 5).Value = "  +0.07%  "

$ws.Cells.Item(8, 5).Value = "  -2.30%  "

$ws.Cells.Item(9, 2).Value = "LidoStakedEther"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Cells.Item(9, 4).Value = "2.569.21"
$ws.Cells.Item(9, 5).Value = "  -4.28%  "

$ws.Cells.Item(10, 2).Value = "Dogecoin"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(10, 4).Value = "0.167"
$ws.Cells.Item(10, 5).Value = "  +0.66%  "

$ws.Cells.Item(11, 5).Value = "  -0.01%  "

$ws.Cells.Item(12, 4).Value = "0.348"
$ws.Cells.Item(12, 5).Value = "  -1.69%  "

$ws.Cells.Item(13, 4).Value = "4.85"
$ws.Cells.Item(13, 5).Value = "  -2.66%  "

$ws.Cells.Item(14, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(14, 4).Value = "3.041.59"
$ws.Cells.Item(14, 5).Value = "  -4.24%  "

$ws.Cells.Item(15, 2).Value = "ShibaInu"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(15, 4).Value = "0.0000184"
$ws.Cells.Item(15, 5).Value = "  +0.09%  "

$ws.Cells.Item(16, 4).Value = "70.646.26"
$ws.Cells.Item(16, 5).Value = "  -1.66%  "

$ws.Cells.Item(17, 4).Value = "25.28"
$ws.Cells.Item(17, 5).Value = "  -3.28%  "

$ws.Cells.Item(18, 4).Value = "2.575.54"
$ws.Cells.Item(18, 5).Value = "  -4.06%  "

$ws.Cells.Item(19, 4).Value = "11.82"
$ws.Cells.Item(19, 5).Value = "  -3.32%  "

$ws.Cells.Item(20, 4).Value = "364.10"
$ws.Cells.Item(20, 5).Value = "  -1.72%  "

$ws.Cells.Item(21, 4).Value = "7.44"
$ws.Cells.Item(21, 5).Value = "  -8.40%  "

$ws.Cells.Item(22, 5).Value = "  -4.51%  "

$ws.Cells.Item(23, 5).Value = "  +1.47%  "

$ws.Cells.Item(24, 5).Value = "  +0.13%  "

$ws.Cells.Item(25, 4).Value = "70.22"
$ws.Cells.Item(25, 5).Value = "  -2.71%  "

$ws.Cells.Item(26, 5).Value = "  -4.17%  "

$ws.Cells.Item(27, 4).Value = "9.36"
$ws.Cells.Item(27, 5).Value = "  -4.04%  "

$ws.Cells.Item(29, 5).Value = "  +0.19%  "

$ws.Cells.Item(30, 4).Value = "0.0₃0934"
$ws.Cells.Item(30, 5).Value = "  -2.53%  "

$ws.Cells.Item(31, 4).Value = "7.84"
$ws.Cells.Item(31, 5).Value = "  -2.70%  "

$ws.Cells.Item(32, 4).Value = "487.02"
$ws.Cells.Item(32, 5).Value = "  -2.20%  "

$ws.Cells.Item(33, 5).Value = "  +1.31%  "

$ws.Cells.Item(34, 5).Value = "  -2.45%  "

$ws.Cells.Item(35, 4).Value = "1.00"
$ws.Cells.Item(35, 5).Value = "  +0.05%  "

$ws.Cells.Item(36, 4).Value = "157.87"
$ws.Cells.Item(36, 5).Value = "  -3.68%  "

$ws.Cells.Item(37, 5).Value = "  +6.32%  "

$ws.Cells.Item(38, 4).Value = "18.80"
$ws.Cells.Item(38, 5).Value = "  -3.95%  "

$ws.Cells.Item(39, 5).Value = "  -1.26%  "

$ws.Cells.Item(40, 5).Value = "  -2.84%  "

$ws.Cells.Item(42, 4).Value = "1.69"
$ws.Cells.Item(42, 5).Value = "  -4.88%  "

$ws.Cells.Item(43, 4).Value = "2.49"
$ws.Cells.Item(43, 5).Value = "  -2.00%  "

$ws.Cells.Item(44, 4).Value = "4.79"
$ws.Cells.Item(44, 5).Value = "  -4.08%  "

$ws.Cells.Item(45, 5).Value = "  -3.52%  "

$ws.Cells.Item(46, 4).Value = "38.53"
$ws.Cells.Item(46, 5).Value = "  -2.05%  "

$ws.Cells.Item(47, 4).Value = "146.42"
$ws.Cells.Item(47, 5).Value = "  -6.99%  "

$ws.Cells.Item(48, 4).Value = "3.58"
$ws.Cells.Item(48, 5).Value = "  -3.92%  "

$ws.Cells.Item(49, 4).Value = "0.535"
$ws.Cells.Item(49, 5).Value = "  -5.23%  "

$ws.Cells.Item(50, 5).Value = "  -5.92%  "

$ws.Cells.Item(51, 4).Value = "0.598"
$ws.Cells.Item(51, 5).Value = "  -1.07%  "
